# Insert a new data row at row 147 (shifts existing rows 147..273 down to 148..274)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(147).Insert()

$ws.Cells.Item(147, 1).Value = 7
$ws.Cells.Item(147, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(147, 3).Value = "Ñuble"
$ws.Cells.Item(147, 4).Value = 44827
$ws.Cells.Item(147, 5).Value = 16
$ws.Cells.Item(147, 6).Value = 100112003
$ws.Cells.Item(147, 7).Value = "Ajo"
$ws.Cells.Item(147, 8).Value = "Chino"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 40
$ws.Cells.Item(147, 11).Value = 23000
$ws.Cells.Item(147, 12).Value = 24000
$ws.Cells.Item(147, 13).Value = 23500
$ws.Cells.Item(147, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(147, 15).Value = "China"
$ws.Cells.Item(147, 16).Value = 2350
$ws.Cells.Item(147, 17).Value = 10
$ws.Cells.Item(147, 18).Value = "Hortaliza"
